$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells I1 ("I0") and J1 ("IF"), reusing the same formatting
# (bold font, border, centered/top alignment) already used by the other
# header cells in row 1.
$ws.Range("H1").Copy()
$ws.Range("I1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("I1").Value = "I0"

$ws.Range("H1").Copy()
$ws.Range("J1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("J1").Value = "IF"

# Fill data rows 2..19:
#   I column = 1 (constant)
#   J column = copy of H column's value
for ($r = 2; $r -le 19; $r++) {
    $ws.Cells.Item($r, 9).Value = 1
    $ws.Cells.Item($r, 10).Value = $ws.Cells.Item($r, 8).Value2
}
